$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scan-Based Contrastive")

# --- New data block in columns E:G, mirroring A:C (Iteration / Dice Score / Hausdorff Distance) ---

# Header row
$ws.Range("E1").Value = "Iteration"
$ws.Range("F1").Value = "Dice Score"
$ws.Range("G1").Value = "Hausdorff Distance"

# Data rows 2-11
$iterData = @(
  @(1, 0.4336, 18.3917),
  @(2, 0.4362, 31.4823),
  @(3, 0.4164, 32.6523),
  @(4, 0.4263, 32.2428),
  @(5, 0.8206, 27.7106),
  @(6, 0.4080, 40.3763),
  @(7, 0.4342, 31.8419),
  @(8, 0.4330, 25.5872),
  @(9, 0.8727, 16.2973),
  @(10, 0.4185, 22.0409)
)

for ($i = 0; $i -lt $iterData.Length; $i++) {
  $r = 2 + $i
  $ws.Range("E$r").Value = $iterData[$i][0]
  $ws.Range("F$r").Value = $iterData[$i][1]
  $ws.Range("G$r").Value = $iterData[$i][2]
}

# Average row
$ws.Range("E12").Value = "Average"
$ws.Range("F12").Formula = "=AVERAGE(F2:F11)"
$ws.Range("G12").Formula = "=AVERAGE(G2:G11)"

# Standard deviation row
$ws.Range("E13").Value = "Standard Deviation"
$ws.Range("F13").Formula = "=_xlfn.STDEV.S(F2:F11)"
$ws.Range("G13").Formula = "=_xlfn.STDEV.S(G2:G11)"

# Copy formatting from A1:C13 onto E1:G13 so the new block matches the existing look
$ws.Range("A1:C13").Copy()
$ws.Range("E1:G13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New column widths for E, F, G
$ws.Range("E1").ColumnWidth = 16.8
$ws.Range("F1").ColumnWidth = 16.1
$ws.Range("G1").ColumnWidth = 17.6

# Labels describing the two data blocks
$ws.Range("A15").Value = "Full Data"
$ws.Range("E15").Value = "Reduced Data (45% of total)"

# Select H11 and make this the active sheet/tab
$ws.Activate()
$ws.Range("H11").Select()
